$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sort fix: Israel moves up in rank (was row 38, now row 36) ---
# Shift country names for rows 36-38: Israel, Ucrania, Bielorrusia
$ws.Range("A36").Value = "Israel"
$ws.Range("A37").Value = "Ucrania"
$ws.Range("A38").Value = "Bielorrusia"

# --- Update "last updated" timestamp in A1 ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Julio de 2020 a las 19:31"

# --- Update statistics data cells (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---
# Row 4
$ws.Range("B4").Value = 4522140
$ws.Range("C4").Value = 23797
$ws.Range("D4").Value = 2208438
$ws.Range("E4").Value = 2160849
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 512
$ws.Range("H4").Value = 152853

# Row 5
$ws.Range("B5").Value = 2498668
$ws.Range("C5").Value = 14019
$ws.Range("D5").Value = 1721560
$ws.Range("E5").Value = 688316
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 158
$ws.Range("H5").Value = 88792

# Row 6
$ws.Range("B6").Value = 1584219
$ws.Range("C6").Value = 52084
$ws.Range("D6").Value = 1020337
$ws.Range("E6").Value = 528882
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 776
$ws.Range("H6").Value = 35000

# Row 12
$ws.Range("B12").Value = 329721
$ws.Range("C12").Value = 2031
$ws.Range("D12").Value = 0
$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 28441

# Row 20
$ws.Range("B20").Value = 228924
$ws.Range("C20").Value = 942
$ws.Range("D20").Value = 212557
$ws.Range("E20").Value = 10708
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 5659

# Row 23
$ws.Range("D23").Value = 77855
$ws.Range("E23").Value = 92300
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 21
$ws.Range("H23").Value = 3200

# Row 32
$ws.Range("B32").Value = 83193
$ws.Range("C32").Value = 914
$ws.Range("D32").Value = 35572
$ws.Range("E32").Value = 41998
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 39
$ws.Range("H32").Value = 5623

# Row 36
$ws.Range("B36").Value = 67734
$ws.Range("C36").Value = 1441
$ws.Range("D36").Value = 32722
$ws.Range("E36").Value = 34521
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 491

# Row 37
$ws.Range("B37").Value = 67597
$ws.Range("C37").Value = 1022
$ws.Range("D37").Value = 37394
$ws.Range("E37").Value = 28553
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 21
$ws.Range("H37").Value = 1650

# Row 38
$ws.Range("B38").Value = 67518
$ws.Range("C38").Value = 152
$ws.Range("D38").Value = 61442
$ws.Range("E38").Value = 5528
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 548

# Row 46
$ws.Range("B46").Value = 50613
$ws.Range("C46").Value = 203
$ws.Range("D46").Value = 35875
$ws.Range("E46").Value = 13013
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 3
$ws.Range("H46").Value = 1725

# Row 60
$ws.Range("B60").Value = 29229
$ws.Range("C60").Value = 614
$ws.Range("D60").Value = 19592
$ws.Range("E60").Value = 8451
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 12
$ws.Range("H60").Value = 1186

# Row 61
$ws.Range("B61").Value = 25942
$ws.Range("C61").Value = 13
$ws.Range("D61").Value = 23364
$ws.Range("E61").Value = 814
$ws.Range("F61").Value = 0

# Row 65
$ws.Range("B65").Value = 22213
$ws.Range("C65").Value = 826
$ws.Range("D65").Value = 17125
$ws.Range("E65").Value = 4754
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 7
$ws.Range("H65").Value = 334

# Row 85
$ws.Range("B85").Value = 9961
$ws.Range("C85").Value = 156
$ws.Range("D85").Value = 6655
$ws.Range("E85").Value = 3106
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 200

# Row 89
$ws.Range("B89").Value = 7647
$ws.Range("C89").Value = 85
$ws.Range("D89").Value = 6176
$ws.Range("E89").Value = 1428
$ws.Range("F89").Value = 0

# Row 95
$ws.Range("B95").Value = 6533
$ws.Range("C95").Value = 158
$ws.Range("D95").Value = 4959
$ws.Range("E95").Value = 1460
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 114

# Row 104
$ws.Range("B104").Value = 4336
$ws.Range("C104").Value = 57
$ws.Range("D104").Value = 1374
$ws.Range("E104").Value = 2759
$ws.Range("F104").Value = 0

# Row 105
$ws.Range("B105").Value = 4205
$ws.Range("C105").Value = 182
$ws.Range("D105").Value = 1753
$ws.Range("E105").Value = 2397
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 55

# Row 132
$ws.Range("B132").Value = 1803
$ws.Range("C132").Value = 17
$ws.Range("D132").Value = 1355
$ws.Range("E132").Value = 381
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 1
$ws.Range("H132").Value = 67

# Row 138
$ws.Range("B138").Value = 1488
$ws.Range("C138").Value = 20
$ws.Range("D138").Value = 1178
$ws.Range("E138").Value = 260
$ws.Range("F138").Value = 0

# Row 141
$ws.Range("B141").Value = 1187
$ws.Range("C141").Value = 5
$ws.Range("D141").Value = 1049
$ws.Range("E141").Value = 127
$ws.Range("F141").Value = 0

# Row 143
$ws.Range("E143").Value = 209
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 17

# Row 147
$ws.Range("B147").Value = 1080
$ws.Range("C147").Value = 13
$ws.Range("D147").Value = 852
$ws.Range("E147").Value = 209
$ws.Range("F147").Value = 0

# Row 149
$ws.Range("D149").Value = 813
$ws.Range("E149").Value = 38
$ws.Range("F149").Value = 0

# Row 150
$ws.Range("B150").Value = 918
$ws.Range("C150").Value = 11
$ws.Range("D150").Value = 804
$ws.Range("E150").Value = 62
$ws.Range("F150").Value = 0

# Row 190
$ws.Range("B190").Value = 88
$ws.Range("C190").Value = 1
$ws.Range("D190").Value = 85
$ws.Range("E190").Value = 2
$ws.Range("F190").Value = 0
